$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 74, shifting existing rows 74-83 down to 75-84
$ws.Rows("74:74").Insert()

# Populate the newly inserted row 74 with the new weekly price record
$ws.Range("A74").Value = 10
$ws.Range("B74").Value = "Vega Modelo de Temuco"
$ws.Range("C74").Value = "La Araucanía"
$ws.Range("D74").Value = 44476
$ws.Range("E74").Value = 9
$ws.Range("F74").Value = 100112012
$ws.Range("G74").Value = "Espinaca"
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 40
$ws.Range("K74").Value = 11000
$ws.Range("L74").Value = 11000
$ws.Range("M74").Value = 11000
$ws.Range("N74").Value = "`$/docena de atados"
$ws.Range("O74").Value = "Región de La Araucanía"
$ws.Range("P74").Value = 3667
$ws.Range("Q74").Value = 3
$ws.Range("R74").Value = "Hortaliza"

# Ensure the date cell keeps the same date number format style used by the other rows (s="2")
$ws.Range("D74").NumberFormat = $ws.Range("D75").NumberFormat
